$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Remember which sheet/cell was active so we can restore the selection
# afterwards -- selecting a cell on "settings" below would otherwise make
# it the active sheet, which is not part of this edit.
$originalActiveSheet = $wb.ActiveSheet.Name

# The "form_id" column (column B) is being dropped from the settings sheet.
# Comments are anchored to a fixed cell and do NOT slide over when a column
# is deleted, so first line up each remaining header's comment with the
# column it will occupy once everything shifts left, and remove the
# trailing comment that belonged to the rightmost (now nonexistent) column.

$ws.Range("B1").Comment.Text("The unique version code that identifies the current state of the form. A common convention is to use a format like yyyymmddrr. For example, 2017021501 is the 1st revision from Feb 15th, 2017.

By default, this template uses a formula to create a date-based version that will update automatically.")

$ws.Range("C1").Comment.Text('Set to ‘pages’ to indicate that groups with the `field-list` appearance represent separate form pages (and all other questions will be shown on their own page). ')

$ws.Range("D1").Comment.Text('Custom namespaces supported in the form.  `cht` must be included here to use the custom `instance::cht` columns on the survey sheet.')

$ws.Range("E1").Comment.Delete()

# Now remove the form_id column itself; version/style/namespaces shift one
# column to the left (C->B, D->C, E->D).
$ws.Columns("B").Delete()

# Move the active-cell marker on the settings sheet the same way Excel
# would leave it after a column delete near the top of the sheet.
$ws.Range("A6").Select()

# Restore whichever sheet was active before we touched "settings".
$wb.Worksheets.Item($originalActiveSheet).Activate()
